$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising | Quicksilver
$ws.Range("L2").Value = 242.42857
$ws.Range("N2").Value = -468.42857
$ws.Range("M2").Value = 25.833336
$ws.Range("H2").Value = 170.76923
$ws.Range("J2").Value = 242.42857
$ws.Range("K2").Value = 87.166664
$ws.Range("I2").Value = 87.166664
# Row 28: The Writing Is Not on the Wall | Enchanted Silver Ink
$ws.Range("L28").Value = 2627.7144
$ws.Range("N28").Value = -3597.7144
$ws.Range("M28").Value = -1373.0769
$ws.Range("H28").Value = 2127.45
$ws.Range("J28").Value = 2627.7144
$ws.Range("K28").Value = 1858.0769
$ws.Range("I28").Value = 1858.0769
# Row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws.Range("I62").Value = 16715.555
$ws.Range("M62").Value = -16091.555
$ws.Range("K62").Value = 16715.555
$ws.Range("H62").Value = 305200.25
# Row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws.Range("I65").Value = 16715.555
$ws.Range("M65").Value = -80457.77499999999
$ws.Range("K65").Value = 83577.77499999999
$ws.Range("H65").Value = 305200.25
# Row 92: Whinier than the Sword | Enchanted Koppranickel Ink
$ws.Range("I92").Value = 26.6
$ws.Range("M92").Value = 1221.4
$ws.Range("K92").Value = 26.6
$ws.Range("H92").Value = 47.57143
# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("I100").Value = 2663
$ws.Range("M100").Value = -2122
$ws.Range("K100").Value = 2663
$ws.Range("H100").Value = 3071.3333
# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("L137").Value = 4841409.300000001
$ws.Range("N137").Value = -4846509.300000001
$ws.Range("H137").Value = 727284.1
$ws.Range("J137").Value = 1613803.1

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("I2").Value = 965.05554
$ws.Range("M2").Value = -852.05554
$ws.Range("K2").Value = 965.05554
$ws.Range("H2").Value = 1099.3077
# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("I45").Value = 3449.111
$ws.Range("M45").Value = -3072.111
$ws.Range("K45").Value = 3449.111
$ws.Range("H45").Value = 3485.6365
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("I61").Value = 1893
$ws.Range("M61").Value = -1681
$ws.Range("K61").Value = 1893
$ws.Range("H61").Value = 79087.38
# Row 76: Sometimes the South Wins | Titanium Mail of Fending
$ws.Range("L76").Value = 121423.57
$ws.Range("N76").Value = -122099.57
$ws.Range("H76").Value = 121423.57
$ws.Range("J76").Value = 121423.57
# Row 79: The Thriller of Autumn (L) | Titanium Mail of Fending
$ws.Range("L79").Value = 121423.57
$ws.Range("N79").Value = -123763.57
$ws.Range("H79").Value = 121423.57
$ws.Range("J79").Value = 121423.57
# Row 116: No Scope | Titanbronze Ingot
$ws.Range("I116").Value = 965.05554
$ws.Range("M116").Value = 1328.94446
$ws.Range("K116").Value = 965.05554
$ws.Range("H116").Value = 1099.3077
# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("I136").Value = 1893
$ws.Range("M136").Value = -3129
$ws.Range("K136").Value = 5679
$ws.Range("H136").Value = 79087.38
# Row 138: Don't Ask about the Rivets | Titanium Gold Helm of Casting
$ws.Range("L138").Value = 94996
$ws.Range("N138").Value = -105276
$ws.Range("H138").Value = 94996
$ws.Range("J138").Value = 94996

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("I3").Value = 965.05554
$ws.Range("M3").Value = -851.05554
$ws.Range("K3").Value = 965.05554
$ws.Range("H3").Value = 1099.3077
# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("L20").Value = 2317.6
$ws.Range("N20").Value = -2811.6
$ws.Range("M20").Value = -1625.1364
$ws.Range("H20").Value = 2011.3438
$ws.Range("J20").Value = 2317.6
$ws.Range("K20").Value = 1872.1364
$ws.Range("I20").Value = 1872.1364
# Row 94: High Steal | High Steel Nugget
$ws.Range("I94").Value = 1676.75
$ws.Range("M94").Value = -1225.75
$ws.Range("K94").Value = 1676.75
$ws.Range("H94").Value = 1743.8572
# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("I99").Value = 2231.8
$ws.Range("M99").Value = -733.8000000000002
$ws.Range("K99").Value = 2231.8
$ws.Range("H99").Value = 2607002.5
# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("I105").Value = 72910.64
$ws.Range("M105").Value = -71163.64
$ws.Range("K105").Value = 72910.64
$ws.Range("H105").Value = 57832.723
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("L134").Value = 11498.1819
$ws.Range("N134").Value = -16568.1819
$ws.Range("M134").Value = -1417.3752
$ws.Range("H134").Value = 2107.9714
$ws.Range("J134").Value = 3832.7273
$ws.Range("K134").Value = 3952.3752
$ws.Range("I134").Value = 1317.4584

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 28: Militia on My Mind | Iron Lance
$ws.Range("L28").Value = 13623.75
$ws.Range("N28").Value = -14113.75
$ws.Range("H28").Value = 13623.75
$ws.Range("J28").Value = 13623.75
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("I31").Value = 4317.75
$ws.Range("M31").Value = -4022.75
$ws.Range("K31").Value = 4317.75
$ws.Range("H31").Value = 5554.1577
# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("I34").Value = 4317.75
$ws.Range("M34").Value = -4115.75
$ws.Range("K34").Value = 4317.75
$ws.Range("H34").Value = 5554.1577
# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("I105").Value = 1090
$ws.Range("M105").Value = 657
$ws.Range("K105").Value = 1090
$ws.Range("H105").Value = 3139.2307
# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("L122").Value = 12975
$ws.Range("N122").Value = -17875
$ws.Range("M122").Value = -2832.1819
$ws.Range("H122").Value = 2444.5334
$ws.Range("J122").Value = 4325
$ws.Range("K122").Value = 5282.1819
$ws.Range("I122").Value = 1760.7273

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 31: Food Fight | Shepherd's Pie
$ws.Range("I31").Value = 2474.6667
$ws.Range("M31").Value = -7136.000100000001
$ws.Range("K31").Value = 7424.000100000001
$ws.Range("H31").Value = 2474.6667
# Row 74: The Nutcracker's Sweets | Royal Eggs
$ws.Range("I74").Value = 2221
$ws.Range("M74").Value = -5602
$ws.Range("K74").Value = 6663
$ws.Range("H74").Value = 2221
# Row 76: Old Victories, New Tastes | Dhalmel Fricassee
$ws.Range("I76").Value = 5971
$ws.Range("M76").Value = -17530
$ws.Range("K76").Value = 17913
$ws.Range("H76").Value = 12728.5
# Row 77: Time for a Midnight Snack (L) | Royal Eggs
$ws.Range("I77").Value = 2221
$ws.Range("M77").Value = -14685
$ws.Range("K77").Value = 19989
$ws.Range("H77").Value = 2221
# Row 79: The Eats of Authenticity (L) | Dhalmel Fricassee
$ws.Range("I79").Value = 5971
$ws.Range("M79").Value = -16587
$ws.Range("K79").Value = 17913
$ws.Range("H79").Value = 12728.5
# Row 82: Persuasion of a Higher Power | Baked Pipira Pira
$ws.Range("I82").Value = 8004.6665
$ws.Range("M82").Value = -23607.9995
$ws.Range("K82").Value = 24013.9995
$ws.Range("H82").Value = 8004.6665
# Row 85: Loaves and Fishes (L) | Baked Pipira Pira
$ws.Range("I85").Value = 8004.6665
$ws.Range("M85").Value = -22609.9995
$ws.Range("K85").Value = 24013.9995
$ws.Range("H85").Value = 8004.6665

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("L126").Value = 12056.25
$ws.Range("N126").Value = -16996.25
$ws.Range("M126").Value = -4069.400000000001
$ws.Range("H126").Value = 3311.4614
$ws.Range("J126").Value = 4018.75
$ws.Range("K126").Value = 6539.400000000001
$ws.Range("I126").Value = 2179.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("I7").Value = 32354.857
$ws.Range("M7").Value = -32242.857
$ws.Range("K7").Value = 32354.857
$ws.Range("H7").Value = 50837.4
# Row 46: Supply Side Logic | Boar Leather
$ws.Range("I46").Value = 1463.3334
$ws.Range("M46").Value = -1275.3334
$ws.Range("K46").Value = 1463.3334
$ws.Range("H46").Value = 2226
# Row 93: Hide to Go Seek | Gagana Leather
$ws.Range("L93").Value = 1611.2
$ws.Range("N93").Value = -4107.2
$ws.Range("M93").Value = -185.3334
$ws.Range("H93").Value = 1544.5
$ws.Range("J93").Value = 1611.2
$ws.Range("K93").Value = 1433.3334
$ws.Range("I93").Value = 1433.3334
# Row 104: Brace Yourselves | Gazelleskin Bracers of Fending
$ws.Range("L104").Value = 10199.2
$ws.Range("N104").Value = -17187.2
$ws.Range("H104").Value = 10199.2
$ws.Range("J104").Value = 10199.2
# Row 122: Hell on Leather | Gaja Leather
$ws.Range("I122").Value = 81582.2
$ws.Range("M122").Value = -242296.6
$ws.Range("K122").Value = 244746.6
$ws.Range("H122").Value = 25052364
# Row 126: Battered Books | Saiga Leather
$ws.Range("I126").Value = 32354.857
$ws.Range("M126").Value = -94594.571
$ws.Range("K126").Value = 97064.571
$ws.Range("H126").Value = 50837.4
# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("L132").Value = 0
$ws.Range("H132").Value = 35850
$ws.Range("J132").Value = 0
$ws.Range("N132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("I107").Value = 6240.1816
$ws.Range("M107").Value = -16800.5448
$ws.Range("K107").Value = 18720.5448
$ws.Range("H107").Value = 8391.8125
# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("L122").Value = 14100
$ws.Range("N122").Value = -19000
$ws.Range("M122").Value = -6903.700000000001
$ws.Range("H122").Value = 3483
$ws.Range("J122").Value = 4700
$ws.Range("K122").Value = 9353.700000000001
$ws.Range("I122").Value = 3117.9
# Row 125: Color Coated | Almasty Serge Coat of Healing
$ws.Range("L125").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("N125").ClearContents()
# Row 126: A Polished Purchase | Snow Linen
$ws.Range("L126").Value = 6750
$ws.Range("N126").Value = -11690
$ws.Range("M126").Value = -682.625
$ws.Range("H126").Value = 1184.1111
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 3152.625
$ws.Range("I126").Value = 1050.875
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("I132").Value = 869.10345
$ws.Range("M132").Value = -77.31034999999974
$ws.Range("K132").Value = 2607.31035
$ws.Range("H132").Value = 1279787.6
# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1499.25
$ws.Range("H136").Value = 1349.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4049.25
$ws.Range("I136").Value = 1349.75
$ws.Range("N136").ClearContents()
